$wb = $excel.ActiveWorkbook

# --- 1. DictionaryPage: remove the 4 rows that duplicated Term/Search Dictionary
#        entries, leaving only the Drug/Genetics Dictionary rows ---
$dp = $wb.Worksheets.Item("DictionaryPage")
$dp.Range("A2:B5").EntireRow.Delete()
$dp.Columns.Item(1).ColumnWidth = 42.33

# --- 2. Duplicate the TermsSpanish sheet to create the new TermsEnglish sheet,
#        inserted immediately before TermsSpanish (same spot TermsSpanish used
#        to occupy) ---
$esOriginal = $wb.Worksheets.Item("TermsSpanish")
$esOriginal.Copy($esOriginal)

# The copy is placed right before "TermsSpanish" and auto-named "TermsSpanish (2)".
# Re-fetch both sheets by name (a variable captured before Copy() tracks the
# new copy afterwards, not the original) and rename the copy.
$enTerms = $wb.Worksheets.Item("TermsSpanish (2)")
$enTerms.Name = "TermsEnglish"

$enTerms = $wb.Worksheets.Item("TermsEnglish")
$esTerms = $wb.Worksheets.Item("TermsSpanish")

# --- 3. Add the new "Dictionary Contains" row to TermsSpanish ---
$esTerms.Range("A5").Value = "/espanol/publicaciones/diccionario/buscar?contains=true&q=tumor"
$esTerms.Range("B5").Value = "Dictionary Contains"
$esTerms.Columns.Item(1).ColumnWidth = 62.15
$esTerms.Range("A6").Select() | Out-Null

# --- 4. Populate TermsEnglish with the English cancer-terms dictionary rows ---
$enTerms.Range("A2").Value = "/publications/dictionaries/cancer-terms"
$enTerms.Range("B2").Value = "Dictionary Page"
$enTerms.Range("A3").Value = "/publications/dictionaries/cancer-terms/search"
$enTerms.Range("B3").Value = "Dictionary Search Page"
$enTerms.Range("A4").Value = "/publications/dictionaries/cancer-terms?expand=D"
$enTerms.Range("B4").Value = "Dictionary Page Expand"
$enTerms.Range("A5").Value = "/publications/dictionaries/cancer-terms/search?contains=true&q=breast"
$enTerms.Range("B5").Value = "Dictionary Contains"
$enTerms.Columns.Item(1).ColumnWidth = 66.5
$enTerms.Range("A3").Select() | Out-Null

# --- 5. DictionaryPage becomes the active tab/selection ---
$dp.Select() | Out-Null
$dp.Range("A4").Select() | Out-Null
